# Update Name of Algo
# Apply updated values to column A for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -22.41730000000002
    10 = -22.0963
    12 = -21.418
    18 = -22.17430000000001
    37 = -19.92069999999999
    55 = -22.183
    68 = -21.4656
    77 = -20.28899999999999
    78 = -19.67679999999998
    81 = -22.11570000000001
    82 = -21.67280000000001
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
